$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column before column C. This shifts every column
#    from C onward one position to the right (C->D, D->E, ... AX->AY), and
#    updates the sheet's used range / merged cells / column widths accordingly.
$ws.Columns("C").Insert()

# 2. The cell that used to be D2:D3 ("Don vi quan ly") now lives at E2:E3
#    after the shift above. Copy its formatting (fill, borders, font,
#    alignment) onto the freshly inserted C2:C3 so the new header column
#    matches the existing merged-header look used elsewhere in row 2-3.
$ws.Range("E2:E3").Copy()
$ws.Range("C2:C3").PasteSpecial(-4122)  # xlPasteFormats

# 3. Put the new header text in the top cell of the pair, then merge the
#    pair together, mirroring how the other header cells (A2:A3, B2:B3,
#    D2:D3, etc.) are merged across the two header rows.
$ws.Range("C2").Value = "Mã quy hoạch vị trí"
$ws.Range("C2:C3").Merge()

# 4. The worksheet also carries a pile of tiny invisible text-box shapes
#    that were all anchored at (column J, row 4) before the edit. Because
#    a new column was inserted before them, they need to move one column
#    to the right (to column K) to keep sitting over the same data column.
#    One of them originally carried a small sub-column offset (28575 EMU
#    = 2.25 pt); the rest sit flush on the column boundary.
$newColLeft = $ws.Cells.Item(4, 11).Left   # left edge of column K (1-indexed 11)
$smallOffsetPts = 28575 / 12700.0

for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)
    if ($i -eq 22) {
        $shp.Left = $newColLeft + $smallOffsetPts
    } else {
        $shp.Left = $newColLeft
    }
}

# 5. Cosmetic: match the saved selection/active cell recorded in the file.
$ws.Range("C7").Select()
